$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly data block (above the old
# row 9), pushing the existing rows 9-14 down to rows 11-16.
$ws.Rows("9:10").Insert()

# New week's data: row 9
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Terminal La Palmera de La Serena"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 44790
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 100112013
$ws.Range("G9").Value = "Alcachofa"
$ws.Range("H9").Value = "Española"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 400
$ws.Range("K9").Value = 14000
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = 14500
$ws.Range("N9").Value = "$/caja 30 unidades"
$ws.Range("O9").Value = "Provincia de Limarí"
$ws.Range("P9").Value = 483
$ws.Range("Q9").Value = 30
$ws.Range("R9").Value = "Hortaliza"

# New week's data: row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Terminal La Palmera de La Serena"
$ws.Range("C10").Value = "Coquimbo"
$ws.Range("D10").Value = 44790
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 100112013
$ws.Range("G10").Value = "Alcachofa"
$ws.Range("H10").Value = "Madrigal"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 400
$ws.Range("K10").Value = 11500
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 11750
$ws.Range("N10").Value = "$/caja 40 unidades"
$ws.Range("O10").Value = "Provincia del Elquí"
$ws.Range("P10").Value = 294
$ws.Range("Q10").Value = 40
$ws.Range("R10").Value = "Hortaliza"
